$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Data")

# Row 2: pressure_1 average
$ws.Range("J2").Value = "1 Avg."
$ws.Range("K2").Formula = "=AVERAGE(D2:D42)"

# Row 3: pressure_2 average
$ws.Range("J3").Value = "2 Avg"
$ws.Range("K3").Formula = "=AVERAGE(E2:E42)"

# Row 4: pressure_3 - note non-aggregated range reference (as in source)
$ws.Range("J4").Value = "3 Avg"
$ws.Range("K4").Formula = "=F2:F42"

# Row 5: pressure_4 average
$ws.Range("J5").Value = "4 Avg"
$ws.Range("K5").Formula = "=AVERAGE(G2:G42)"

# Row 6: pressure_5 average
$ws.Range("J6").Value = "5 Avg"
$ws.Range("K6").Formula = "=AVERAGE(H2:H42)"

# Update selection to match target (O7)
$ws.Range("O7").Select()
